$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) and "Volumen" (M) values between rows 3 and 4.
$ws.Range("D3").Value = 44253
$ws.Range("D4").Value = 44993

$ws.Range("M3").Value = 12
$ws.Range("M4").Value = 14
